$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns B:E to stay as Text so Excel does not auto-convert
# numeric-looking strings (e.g. "32.48", "1.00") into numbers, matching
# the original inlineStr cell types in the workbook.
$ws.Range("B2:E51").NumberFormat = "@"

# Unicode subscript digits used by a few coin prices (e.g. 0.0 3x / 0.0 6x)
$sub3 = [char]0x2083
$sub6 = [char]0x2086

$ws.Range("D2").Value = '34.427.02'
$ws.Range("E2").Value = '  +0.95%  '
$ws.Range("D3").Value = '1.796.29'
$ws.Range("E3").Value = '  +0.55%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '226.75'
$ws.Range("E5").Value = '  +0.07%  '
$ws.Range("E6").Value = '  +1.42%  '
$ws.Range("E7").Value = '  -0.06%  '
$ws.Range("D8").Value = '32.48'
$ws.Range("E8").Value = '  +1.69%  '
$ws.Range("D9").Value = '0.296'
$ws.Range("E9").Value = '  +1.32%  '
$ws.Range("E10").Value = '  +0.20%  '
$ws.Range("E11").Value = '  +0.62%  '
$ws.Range("D12").Value = '2.054.87'
$ws.Range("E12").Value = '  +0.52%  '
$ws.Range("D13").Value = '11.13'
$ws.Range("E13").Value = '  -0.97%  '
$ws.Range("D14").Value = '1.795.86'
$ws.Range("E14").Value = '  +0.38%  '
$ws.Range("D15").Value = '0.631'
$ws.Range("E15").Value = '  +1.84%  '
$ws.Range("D16").Value = '34.374.39'
$ws.Range("E16").Value = '  +1.01%  '
$ws.Range("E17").Value = '  +0.79%  '
$ws.Range("D19").Value = [string]::Concat('0.0', $sub3, '0802')
$ws.Range("E19").Value = '  +3.10%  '
$ws.Range("D20").Value = '246.69'
$ws.Range("E20").Value = '  +0.38%  '
$ws.Range("D21").Value = '11.01'
$ws.Range("E21").Value = '  +1.77%  '
$ws.Range("E22").Value = '  +0.06%  '
$ws.Range("D23").Value = '4.15'
$ws.Range("E23").Value = '  +1.39%  '
$ws.Range("D24").Value = '2.05'
$ws.Range("E24").Value = '  +0.28%  '
$ws.Range("D25").Value = '163.00'
$ws.Range("E25").Value = '  +0.90%  '
$ws.Range("E26").Value = '  +0.57%  '
$ws.Range("D27").Value = '16.42'
$ws.Range("E27").Value = '  +0.63%  '
$ws.Range("E28").Value = '  +2.11%  '
$ws.Range("E29").Value = '  +0.01%  '
$ws.Range("E30").Value = '  +0.95%  '
$ws.Range("E31").Value = '  +0.19%  '
$ws.Range("D32").Value = '3.91'
$ws.Range("E32").Value = '  +8.26%  '
$ws.Range("D33").Value = '3.79'
$ws.Range("E33").Value = '  +3.61%  '
$ws.Range("E34").Value = '  +1.33%  '
$ws.Range("D35").Value = '1.442.38'
$ws.Range("E35").Value = '  -0.35%  '
$ws.Range("D36").Value = '2.63'
$ws.Range("E36").Value = '  +10.03%  '
$ws.Range("E37").Value = '  +2.82%  '
$ws.Range("E39").Value = '  -0.96%  '
$ws.Range("D40").Value = '83.79'
$ws.Range("E40").Value = '  +4.58%  '
$ws.Range("E41").Value = '  +1.35%  '
$ws.Range("D42").Value = '0.936'
$ws.Range("E42").Value = '  +1.87%  '
$ws.Range("D43").Value = '2.76'
$ws.Range("E43").Value = '  +3.06%  '
$ws.Range("E44").Value = '  +2.35%  '
$ws.Range("E45").Value = '  +3.13%  '
$ws.Range("E46").Value = '  +0.92%  '
$ws.Range("E47").Value = '  +0.02%  '
$ws.Range("D48").Value = '1.951.06'
$ws.Range("E48").Value = '  +0.29%  '
$ws.Range("D49").Value = '105.80'
$ws.Range("E49").Value = '  -1.56%  '
$ws.Range("B50").Value = 'PaxDollar'
$ws.Range("C50").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D50").Value = '1.00'
$ws.Range("E50").Value = '  +0.00%  '
$ws.Range("B51").Value = 'BabyDogeCoin'
$ws.Range("C51").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D51").Value = [string]::Concat('0.0', $sub6, '0129')
$ws.Range("E51").Value = '  -6.04%  '
